# Update the "Latest HO Xliff Generate Date" / handoff / handback timestamps
# for the row corresponding to 0fc0d813-078d-4b8a-9b21-995f72b5ad67.md
# (report regeneration for the handback status workbook).

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
# G4: "Latest HO Xliff Generate Date" -> 2016-08-29 22:48:00
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-29 22:48:00"

# --- zh-cn sheet ---
# H4: Correspond Handoff Datetime  -> 2016-08-29 22:47:55
# K4: Correspond Handback DateTime -> 2016-08-29 22:48:30
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-29 22:47:55"
$wsZhCn.Range("K4").Value = "2016-08-29 22:48:30"

# --- de-de sheet ---
# H4: Correspond Handoff Datetime  -> 2016-08-29 22:48:00
# K4: Correspond Handback DateTime -> 2016-08-29 22:48:37
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-08-29 22:48:00"
$wsDeDe.Range("K4").Value = "2016-08-29 22:48:37"
